{"js": "// The document contains three paragraphs each holding a split\n// \"<id>...</id>\" marker: the literal text \"<id>\" / \"p067v_N\" / \"</id>\"\n// is spread across three separate runs (the <id> and </id> runs are\n// Courier New / color 7f6000, the \"p067v_N\" run is Arial / black).\n// The edit collapses each of those triples into a single run whose\n// text is the full \"<id>p067v_N</id>\" string, keeping the formatting\n// of the first (\"<id>\") run (Courier New, color 7f6000).\n//\n// Doing a Replace-insertText over a Range that already spans all three\n// runs achieves exactly that merge: Word (and this shim) rewrites the\n// matched range as one run using the formatting of the range's first\n// run.\n\nconst body = context.document.body;\n\nfor (let n = 1; n <= 3; n++) {\n  const marker = \"<id>p067v_\" + n + \"</id>\";\n\n  const found = body.search(marker, { matchCase: true, matchWildcards: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(marker, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains three paragraphs each holding a split\n# \"<id>...</id>\" marker: the literal text \"<id>\" / \"p067v_N\" / \"</id>\"\n# is spread across three separate runs (the <id> and </id> runs are\n# Courier New / color 7f6000, the \"p067v_N\" run is Arial / black).\n# The edit collapses each of those triples into a single run whose\n# text is the full \"<id>p067v_N</id>\" string, keeping the formatting\n# of the first (\"<id>\") run (Courier New, color 7f6000).\n#\n# A Find/Replace whose search text is the whole marker (it spans all\n# three runs) performs exactly that merge: Word rewrites the matched\n# range as a single run using the formatting of the first run in the\n# match - i.e. wdReplaceOne (2) with MatchCase on and MatchWildcards\n# off so the angle brackets are treated literally.\n\n$d = $word.ActiveDocument\n\nfor ($n = 1; $n -le 3; $n++) {\n    $marker = \"<id>p067v_$n</id>\"\n\n    $find = $d.Content.Find\n    $find.Text = $marker\n    $find.Replacement.Text = $marker\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #         Format, ReplaceWith, Replace)\n    $find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, $marker, 2)\n}\n"}
